# Update LR-pair data rows (2-11) with new TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Col1a2"
$ws.Cells.Item(2, 3).Value = "Gp6"
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 67.91996
$ws.Cells.Item(2, 8).Value = 203.75988
$ws.Cells.Item(2, 9).Value = 0.02375577759132129
$ws.Cells.Item(2, 10).Value = 0.02375577759132129
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.009511333333333333
$ws.Cells.Item(2, 14).Value = 0.028534
$ws.Cells.Item(2, 15).Value = 0.7967497835981349
$ws.Cells.Item(2, 16).Value = 0.7967497835981349
$ws.Cells.Item(2, 17).Value = 0.6460093795466667
$ws.Cells.Item(2, 18).Value = 5.81408441592
$ws.Cells.Item(2, 19).Value = 0.01892741065509066
$ws.Cells.Item(2, 20).Value = 0.01892741065509066

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Col1a2"
$ws.Cells.Item(3, 3).Value = "Gp6"
$ws.Cells.Item(3, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 67.91996
$ws.Cells.Item(3, 8).Value = 203.75988
$ws.Cells.Item(3, 9).Value = 0.02375577759132129
$ws.Cells.Item(3, 10).Value = 0.02375577759132129
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.002426333333333333
$ws.Cells.Item(3, 14).Value = 0.007279
$ws.Cells.Item(3, 15).Value = 0.2032502164018652
$ws.Cells.Item(3, 16).Value = 0.2032502164018652
$ws.Cells.Item(3, 17).Value = 0.1647964629466667
$ws.Cells.Item(3, 18).Value = 1.48316816652
$ws.Cells.Item(3, 19).Value = 0.004828366936230634
$ws.Cells.Item(3, 20).Value = 0.004828366936230634

# Row 4
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Col1a2"
$ws.Cells.Item(4, 3).Value = "Gp6"
$ws.Cells.Item(4, 4).Value = "FAPs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 2623.51945
$ws.Cells.Item(4, 8).Value = 7870.55835
$ws.Cells.Item(4, 9).Value = 0.9176057312269553
$ws.Cells.Item(4, 10).Value = 0.9176057312269554
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.009511333333333333
$ws.Cells.Item(4, 14).Value = 0.028534
$ws.Cells.Item(4, 15).Value = 0.7967497835981349
$ws.Cells.Item(4, 16).Value = 0.7967497835981349
$ws.Cells.Item(4, 17).Value = 24.95316799543334
$ws.Cells.Item(4, 18).Value = 224.5785119589
$ws.Cells.Item(4, 19).Value = 0.7311021677834849
$ws.Cells.Item(4, 20).Value = 0.731102167783485

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Col1a2"
$ws.Cells.Item(5, 3).Value = "Gp6"
$ws.Cells.Item(5, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 2623.51945
$ws.Cells.Item(5, 8).Value = 7870.55835
$ws.Cells.Item(5, 9).Value = 0.9176057312269553
$ws.Cells.Item(5, 10).Value = 0.9176057312269554
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.002426333333333333
$ws.Cells.Item(5, 14).Value = 0.007279
$ws.Cells.Item(5, 15).Value = 0.2032502164018652
$ws.Cells.Item(5, 16).Value = 0.2032502164018652
$ws.Cells.Item(5, 17).Value = 6.365532692183334
$ws.Cells.Item(5, 18).Value = 57.28979422965
$ws.Cells.Item(5, 19).Value = 0.1865035634434705
$ws.Cells.Item(5, 20).Value = 0.1865035634434705

# Row 6
$ws.Cells.Item(6, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(6, 2).Value = "Col1a2"
$ws.Cells.Item(6, 3).Value = "Gp6"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.376679
$ws.Cells.Item(6, 8).Value = 4.130037
$ws.Cells.Item(6, 9).Value = 0.0004815091195378001
$ws.Cells.Item(6, 10).Value = 0.0004815091195378002
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.009511333333333333
$ws.Cells.Item(6, 14).Value = 0.028534
$ws.Cells.Item(6, 15).Value = 0.7967497835981349
$ws.Cells.Item(6, 16).Value = 0.7967497835981349
$ws.Cells.Item(6, 17).Value = 0.013094052862
$ws.Cells.Item(6, 18).Value = 0.117846475758
$ws.Cells.Item(6, 19).Value = 0.0003836422867922707
$ws.Cells.Item(6, 20).Value = 0.0003836422867922707

# Row 7
$ws.Cells.Item(7, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(7, 2).Value = "Col1a2"
$ws.Cells.Item(7, 3).Value = "Gp6"
$ws.Cells.Item(7, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.376679
$ws.Cells.Item(7, 8).Value = 4.130037
$ws.Cells.Item(7, 9).Value = 0.0004815091195378001
$ws.Cells.Item(7, 10).Value = 0.0004815091195378002
$ws.Cells.Item(7, 11).Value = 1
$ws.Cells.Item(7, 12).Value = 0.3333333333333333
$ws.Cells.Item(7, 13).Value = 0.002426333333333333
$ws.Cells.Item(7, 14).Value = 0.007279
$ws.Cells.Item(7, 15).Value = 0.2032502164018652
$ws.Cells.Item(7, 16).Value = 0.2032502164018652
$ws.Cells.Item(7, 17).Value = 0.003340282147
$ws.Cells.Item(7, 18).Value = 0.030062539323
$ws.Cells.Item(7, 19).Value = 0.00009786683274552948
$ws.Cells.Item(7, 20).Value = 0.00009786683274552948

# Row 8
$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 2).Value = "Col1a2"
$ws.Cells.Item(8, 3).Value = "Gp6"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 164.8447596666666
$ws.Cells.Item(8, 8).Value = 494.534279
$ws.Cells.Item(8, 9).Value = 0.05765632735555414
$ws.Cells.Item(8, 10).Value = 0.05765632735555416
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.009511333333333333
$ws.Cells.Item(8, 14).Value = 0.028534
$ws.Cells.Item(8, 15).Value = 0.7967497835981349
$ws.Cells.Item(8, 16).Value = 0.7967497835981349
$ws.Cells.Item(8, 17).Value = 1.567893457442889
$ws.Cells.Item(8, 18).Value = 14.111041116986
$ws.Cells.Item(8, 19).Value = 0.04593766634360098
$ws.Cells.Item(8, 20).Value = 0.045937666343601

# Row 9
$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 2).Value = "Col1a2"
$ws.Cells.Item(9, 3).Value = "Gp6"
$ws.Cells.Item(9, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 164.8447596666666
$ws.Cells.Item(9, 8).Value = 494.534279
$ws.Cells.Item(9, 9).Value = 0.05765632735555414
$ws.Cells.Item(9, 10).Value = 0.05765632735555416
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.002426333333333333
$ws.Cells.Item(9, 14).Value = 0.007279
$ws.Cells.Item(9, 15).Value = 0.2032502164018652
$ws.Cells.Item(9, 16).Value = 0.2032502164018652
$ws.Cells.Item(9, 17).Value = 0.3999683352045555
$ws.Cells.Item(9, 18).Value = 3.599715016841
$ws.Cells.Item(9, 19).Value = 0.01171866101195316
$ws.Cells.Item(9, 20).Value = 0.01171866101195317

# Row 10
$ws.Cells.Item(10, 1).Value = "Resolving-Mac"
$ws.Cells.Item(10, 2).Value = "Col1a2"
$ws.Cells.Item(10, 3).Value = "Gp6"
$ws.Cells.Item(10, 4).Value = "FAPs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 1.431418
$ws.Cells.Item(10, 8).Value = 4.294254
$ws.Cells.Item(10, 9).Value = 0.0005006547066313635
$ws.Cells.Item(10, 10).Value = 0.0005006547066313636
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.009511333333333333
$ws.Cells.Item(10, 14).Value = 0.028534
$ws.Cells.Item(10, 15).Value = 0.7967497835981349
$ws.Cells.Item(10, 16).Value = 0.7967497835981349
$ws.Cells.Item(10, 17).Value = 0.01361469373733333
$ws.Cells.Item(10, 18).Value = 0.122532243636
$ws.Cells.Item(10, 19).Value = 0.0003988965291659266
$ws.Cells.Item(10, 20).Value = 0.0003988965291659266

# Row 11
$ws.Cells.Item(11, 1).Value = "Resolving-Mac"
$ws.Cells.Item(11, 2).Value = "Col1a2"
$ws.Cells.Item(11, 3).Value = "Gp6"
$ws.Cells.Item(11, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 1.431418
$ws.Cells.Item(11, 8).Value = 4.294254
$ws.Cells.Item(11, 9).Value = 0.0005006547066313635
$ws.Cells.Item(11, 10).Value = 0.0005006547066313636
$ws.Cells.Item(11, 11).Value = 1
$ws.Cells.Item(11, 12).Value = 0.3333333333333333
$ws.Cells.Item(11, 13).Value = 0.002426333333333333
$ws.Cells.Item(11, 14).Value = 0.007279
$ws.Cells.Item(11, 15).Value = 0.2032502164018652
$ws.Cells.Item(11, 16).Value = 0.2032502164018652
$ws.Cells.Item(11, 17).Value = 0.003473097207333333
$ws.Cells.Item(11, 18).Value = 0.03125787486600001
$ws.Cells.Item(11, 19).Value = 0.000101758177465437
$ws.Cells.Item(11, 20).Value = 0.000101758177465437
